$d = $word.ActiveDocument

# Turn on revision tracking before editing: in this engine, a plain
# text mutation inside a paragraph re-coalesces every adjacent run that
# shares identical formatting, which would wreck the paragraph's run
# layout. Tracking the change as an insertion/deletion keeps each
# existing run intact; AcceptAllRevisions() then folds the tracked
# change back into plain runs without re-triggering that coalescing.
$d.TrackRevisions = $true

# ---------------------------------------------------------------------
# Helper: replace exactly one (the next) whole-word, case-sensitive
# match of $findText with $replaceText inside the given story range,
# using a duplicate so the original range keeps covering the full
# story (and repeated calls walk left-to-right through it).
# wdReplaceOne = 1
# ---------------------------------------------------------------------
function Replace-OneMatch($storyRange, $findText, $replaceText) {
    $rng = $storyRange.Duplicate
    $rng.Find.ClearFormatting()
    $rng.Find.Execute($findText, $true, $true, $false, $false, $false, `
                       $true, 1, $false, $replaceText, 1) | Out-Null
}

# ---------------------------------------------------------------------
# 1) Body text: "A TERE," -> "A QWER," (single occurrence in the body)
# ---------------------------------------------------------------------
foreach ($story in $d.StoryRanges) {
    if ($story.StoryType -eq 1) {
        Replace-OneMatch $story "TERE" "QWER"
    }
}

# ---------------------------------------------------------------------
# 2) Header text (wdPrimaryHeaderStory = 7): several runs of
#    TRE / TERE / Tre / tre need to become QWER / QWER / Qwer|Qewr / qwer
#    in document order. Each call only touches the first remaining
#    match, so repeated calls walk left to right through the header.
# ---------------------------------------------------------------------
foreach ($story in $d.StoryRanges) {
    if ($story.StoryType -eq 7) {

        # "DIRETORIA DE ENSINO REGIAO TRE" -> "...QWER"
        Replace-OneMatch $story "TRE" "QWER"

        # "TERE – DEP." -> "QWER – DEP."
        Replace-OneMatch $story "TERE" "QWER"

        # Address line with five "Tre" runs -> Qwer, Qwer, Qewr, Qewr, Qwer
        Replace-OneMatch $story "Tre" "Qwer"
        Replace-OneMatch $story "Tre" "Qwer"
        Replace-OneMatch $story "Tre" "Qewr"
        Replace-OneMatch $story "Tre" "Qewr"
        Replace-OneMatch $story "Tre" "Qwer"

        # CEP / Tel / Email lines, three "tre" runs -> qwer
        Replace-OneMatch $story "tre" "qwer"
        Replace-OneMatch $story "tre" "qwer"
        Replace-OneMatch $story "tre" "qwer"
    }
}

# Fold the tracked insert/delete pairs back into normal runs and turn
# tracking back off.
$d.AcceptAllRevisions() | Out-Null
$d.TrackRevisions = $false
